$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2..11, column B (description) cells ---
# Row 2 keeps its original two rich-text runs; we rebuild full text then
# re-apply character-level formatting so the run split is preserved as closely as possible.

# Row 2
$full2 = 'Manages projects and development teams executing in a range of methodologies including waterfall, agile, and lean; Ensures the project meets scope, schedule, and budget; Serves as the Scrum Master for Agile projects; liaison between business and technical team; performs risk management; ensures government receives actionable information in a timely manner necessary to obtain decision/guidance to facilitate project execution. Minimum Education: A Bachelor’s Degree Minimum Experience: Five (5) years'
$ws.Cells.Item(2, 2).Value2 = $full2
# Row 3
$full3 = 'Manages projects and development teams executing in a range of methodologies including waterfall, agile, and lean; Ensures the project meets scope, schedule, and budget; Serves as the Scrum Master for Agile projects; liaison between business and technical team; performs risk management; ensures government receives actionable information in a timely manner necessary to obtain decision/guidance to facilitate project execution. Minimum Education: A Bachelor’s Degree Minimum Experience: Ten (10) years'
$ws.Cells.Item(3, 2).Value2 = $full3
# Row 4
$full4 = 'Analyzes and defines security requirements for Multilevel Security (MLS) issues. Designs, develops, engineers, and implements solutions to MLS requirements. Responsible for the implementation and development of the MLS. Gathers and organizes technical information about an organization’s mission goals and needs, existing security products, and ongoing programs in the MLS arena. Performs risk analyses, which also include risk assessment. Provides daily supervision and direction to staff. Minimum Education: A Bachelor’s Degree Minimum Experience: Eight (8) years'
$ws.Cells.Item(4, 2).Value2 = $full4
# Row 5
$full5 = 'Analyzes security measures for more than one IT functional area (e.g., data, systems, network and/or Web) across the enterprise. Develops, implements, communicates and provides training of security assessments, policies and procedures Tracks, monitors, and enforces security policies, reviews security violation reports and investigates possible security exceptions, and updates, maintains and documents security controls. Prepares reports on security matters to develop security risk analysis scenarios and response procedures. Evaluates and recommends products and/or procedures to enhance productivity and effectiveness. Minimum Education: A Bachelor’s Degree Minimum Experience: Seven (7) years'
$ws.Cells.Item(5, 2).Value2 = $full5
# Row 6
$full6 = 'Experience with cloud services - including open source technology, software development, system engineering, scripting languages and multiple cloud provider environments. Additionally, Cloud engineers need to be familiar with one or more of the following: OpenStack, Amazon Web Services, Rackspace, Google Compute Engine, Microsoft Azure and Docker. Experience with APIs, orchestration, automation and DevOps are also important. Minimum Education: A Bachelor’s Degree Minimum Experience: Seven (7) years'
$ws.Cells.Item(6, 2).Value2 = $full6
# Row 7
$full7 = 'Data Scientist will have necessary statistical modelling, mathematical, big data analytics and predictive modelling skills to build the required algorithms necessary to ask right questions and build objective visualizations and findings from it. Data Scientist will have knowledge of integrating multiple systems and datasets to provide new insights. Examples of required skillset: • Prior experience working as a data architect and managing information schema for large organizations • Experience with big data analytic tools such as Hadoop, Hive, MapReduce, SPLUNK, Elastic Search • Understanding and good working knowledge of SQL and NoSQL • Experience in machine learning, statistical modelling, and predictive analysis • Extensive experience with a statistical programming language. Minimum Education: A Bachelor’s Degree Minimum Experience: Seven (7) years
'
$ws.Cells.Item(7, 2).Value2 = $full7
# Row 8
$full8 = 'Responsible for creating front-end design solutions for both web and mobile platforms. The role involves working closely with project manager, analyst, developers and testers to determine ideal design solution. Conduct usability testing to make sure design satisfies all project requirements. Required skillset: • Design mock-up templates using a combination of tools such as HTML, CSS, Photoshop and other standard industry design tools. • Develop responsive design in HTML5 and CSS3 for mobile compatibility • Experience with JavaScript • Experience designing graphics and UI for mobile development • Expertise in Adobe Creative Suite • Design custom logos and images • Understanding of up to date web standards and specifications • Experience with distributed source control systems such as git. Minimum Education: A Bachelor’s Degree Minimum Experience: Seven (7) years'
$ws.Cells.Item(8, 2).Value2 = $full8
# Row 9
$full9 = 'Develop, modify, or update applications used by business units or infrastructure units. Lead, or play lead technical role in development teams'' efforts to determine unit needs and business processes that are automated by the application. Participate in or review all of the steps in the software development life cycle to create and modify the software. Minimum Education: A Bachelor’s Degree Minimum Experience: Five (5) Years'
$ws.Cells.Item(9, 2).Value2 = $full9
# Row 10
$full10 = 'Develop, modify, or update applications used by business units or infrastructure units. Lead, or play lead technical role in development teams'' efforts to determine unit needs and business processes that are automated by the application. Participate in or review all of the steps in the software development life cycle to create and modify the software. Minimum Education: A Bachelor’s Degree Minimum Experience: Ten (10) Years'
$ws.Cells.Item(10, 2).Value2 = $full10
# Row 11
$full11 = 'The Test Automation Engineer is responsible for the analysis of project functional requirements as well as development of code in Java for automating test scenarios. The incumbent is responsible for the analysis of functional requirements, testing applications, developing test plans, test cases and test scripts, and evaluating test results to determine compliance with test plans and established business processes. Minimum Education: A Bachelor’s Degree Minimum Experience: Five (5) years'
$ws.Cells.Item(11, 2).Value2 = $full11

# Re-apply the rich-text run split on row 2 (originally two runs; the appended
# clause becomes a third run). Font size 11.5 is not representable through this
# COM host (Font.Size only accepts whole numbers), so the closest achievable
# integer size is used.
$r2run1 = 'Manages projects and development teams executing in a range of methodologies including waterfall, agile, and lean; Ensures the project meets scope, schedule, and budget; Serves as the Scrum Master for Agile projects; liaison between business and technical team; performs risk management; ensures government receives '
$r2run2 = 'actionable information in a timely manner necessary to obtain decision/guidance to facilitate project execution.'
$r2c2 = $ws.Cells.Item(2, 2)
$r2start2 = $r2run1.Length + 1
$r2len2 = $r2run2.Length
$r2start3 = $r2start2 + $r2len2
$r2len3 = $full2.Length - ($r2start3 - 1)
$r2c2.Characters($r2start2, $r2len2).Font.Size = 12
$r2c2.Characters($r2start3, $r2len3).Font.Size = 12

# --- Row heights (rows 2..11) to match the updated (wrapped) text ---
$ws.Rows.Item(2).RowHeight = 119
$ws.Rows.Item(3).RowHeight = 119
$ws.Rows.Item(4).RowHeight = 119
$ws.Rows.Item(5).RowHeight = 153
$ws.Rows.Item(6).RowHeight = 119
$ws.Rows.Item(7).RowHeight = 204
$ws.Rows.Item(8).RowHeight = 187
$ws.Rows.Item(9).RowHeight = 85
$ws.Rows.Item(10).RowHeight = 85
$ws.Rows.Item(11).RowHeight = 102

# --- View state: scroll position + active selection, matching authored file ---
$ws.Range("E11").Select()
try {
  $excel.ActiveWindow.ScrollRow = 8
  $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$ws.Range("E11").Select()

Write-Output "edit complete"
